$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.722.34"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "3.334.15"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.19"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.71"
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +2.08%  "
$ws.Range("D9").Value = "3.330.94"
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("E10").Value = "  +6.37%  "
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.88"
$ws.Range("E12").Value = "  +4.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "691.29"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "3.881.53"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").Value = "67.723.63"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "3.351.54"
$ws.Range("E19").Value = "  +2.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.60"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.04"
$ws.Range("E21").Value = "  +4.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.893"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.44"
$ws.Range("E23").Value = "  +4.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.91"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.59"
$ws.Range("E25").Value = "  +4.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.91"
$ws.Range("E26").Value = "  +2.13%  "
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.46"
$ws.Range("E28").Value = "  +5.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.99"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  +3.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.04"
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "569.21"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.01"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.716.98"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.23"
$ws.Range("E36").Value = "  +3.53%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.06"
$ws.Range("E39").Value = "  +12.16%  "
$ws.Range("E40").Value = "  +4.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.16"
$ws.Range("E41").Value = "  +7.12%  "
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.33"
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0672"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("E45").Value = "  +3.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0408"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.65"
$ws.Range("E47").Value = "  +5.72%  "
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.94"
$ws.Range("E51").Value = "  +3.46%  "
